# Signed Off Time Sheets - As of 28/02/2014
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name (was blank) filled in
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor sign-off row (initials + date), mirroring the employee's
# signature row (A25/D25) directly above the "Supervisor Signature" label
$ws.Range("A27").Value = "P.S"

$signOffDate = Get-Date -Year 2014 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("D27").Value = $signOffDate
$ws.Range("D27").NumberFormat = "m/d/yyyy"

# Move the active selection to the newly completed sign-off cell
$ws.Range("D27:E27").Select()
